# FlowInterface.xlsx: update the "Agent Module" row on the Interface sheet
# so Execute and JurisdictionWiseReport are both switched from "No" to "Yes".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Interface")

$ws.Range("E11").Value = "Yes"
$ws.Range("F11").Value = "Yes"

# Match the author's final cursor position on the sheet.
$ws.Activate()
$ws.Range("E14").Select()
